$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(805, 1).Value = "hp.proxy.yy.duowan.com"
$ws.Cells.Item(805, 2).Value = "多玩游戏网"
$ws.Cells.Item(806, 1).Value = "ls.duowan.com"
$ws.Cells.Item(806, 2).Value = "多玩游戏网"
$ws.Cells.Item(807, 1).Value = "wotbox.duowan.com"
$ws.Cells.Item(807, 2).Value = "多玩游戏网"
$ws.Cells.Item(808, 1).Value = "cname.duowan.com"
$ws.Cells.Item(808, 2).Value = "多玩游戏网"
$ws.Cells.Item(809, 1).Value = "box.dwstatic.com"
$ws.Cells.Item(809, 2).Value = "多玩游戏网"
$ws.Cells.Item(810, 1).Value = "weihuialert.yy.com"
$ws.Cells.Item(810, 2).Value = "多玩游戏网"
$ws.Cells.Item(811, 1).Value = "wuxia.duowan.com"
$ws.Cells.Item(811, 2).Value = "多玩游戏网"
$ws.Cells.Item(812, 1).Value = "mobpush.yy.com"
$ws.Cells.Item(812, 2).Value = "多玩游戏网"
$ws.Cells.Item(813, 1).Value = "thirdlogin.yy.com"
$ws.Cells.Item(813, 2).Value = "多玩游戏网"
$ws.Cells.Item(814, 1).Value = "wot.duowan.com"
$ws.Cells.Item(814, 2).Value = "多玩游戏网"
$ws.Cells.Item(815, 1).Value = "zx.duowan.com"
$ws.Cells.Item(815, 2).Value = "多玩游戏网"
$ws.Cells.Item(816, 1).Value = "kf.yy.com"
$ws.Cells.Item(816, 2).Value = "多玩游戏网"
$ws.Cells.Item(817, 1).Value = "coc.5253.com"
$ws.Cells.Item(817, 2).Value = "多玩游戏网"
$ws.Cells.Item(818, 1).Value = "mc.duowan.com"
$ws.Cells.Item(818, 2).Value = "多玩游戏网"
$ws.Cells.Item(819, 1).Value = "comfrontdl.yy.duowan.com"
$ws.Cells.Item(819, 2).Value = "多玩游戏网"
$ws.Cells.Item(820, 1).Value = "template.dl.yy.com"
$ws.Cells.Item(820, 2).Value = "多玩游戏网"
$ws.Cells.Item(821, 1).Value = "m1.dwstatic.com"
$ws.Cells.Item(821, 2).Value = "多玩游戏网5"
$ws.Cells.Item(822, 1).Value = "earn.yystatic.com"
$ws.Cells.Item(822, 2).Value = "多玩游戏网6"
$ws.Cells.Item(823, 1).Value = "weblbs.yystatic.com"
$ws.Cells.Item(823, 2).Value = "多玩游戏网6"
$ws.Cells.Item(824, 1).Value = "2.dximscreenshot7.yy.yystatic.com"
$ws.Cells.Item(824, 2).Value = "多玩游戏网6"
$ws.Cells.Item(825, 1).Value = "8.dximscreenshot7.yy.yystatic.com"
$ws.Cells.Item(825, 2).Value = "多玩游戏网6"

$ws.Range("F808").Select()
$excel.ActiveWindow.ScrollRow = 792
$excel.ActiveWindow.ScrollColumn = 1
